# Updates market price/profit figures (columns H-N) for specific leve rows
# across all 8 crafting-job sheets, matching the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 215.33333
$ws.Cells.Item(5, 9).Value = 99.8
$ws.Cells.Item(5, 10).Value = 359.75
$ws.Cells.Item(5, 11).Value = 99.8
$ws.Cells.Item(5, 12).Value = 359.75
$ws.Cells.Item(5, 13).Value = 15.2
$ws.Cells.Item(5, 14).Value = -589.75

$ws.Cells.Item(58, 8).Value = 259
$ws.Cells.Item(58, 9).Value = 73.75
$ws.Cells.Item(58, 10).Value = 1000
$ws.Cells.Item(58, 11).Value = 221.25
$ws.Cells.Item(58, 12).Value = 3000
$ws.Cells.Item(58, 13).Value = -71.25
$ws.Cells.Item(58, 14).Value = -3300

$ws.Cells.Item(97, 8).Value = 250000350
$ws.Cells.Item(97, 10).Value = 250000350
$ws.Cells.Item(97, 12).Value = 750001050
$ws.Cells.Item(97, 14).Value = -750002042

$ws.Cells.Item(132, 8).Value = 4595
$ws.Cells.Item(132, 9).Value = 1801.4166
$ws.Cells.Item(132, 10).Value = 9384
$ws.Cells.Item(132, 11).Value = 5404.2498
$ws.Cells.Item(132, 12).Value = 28152
$ws.Cells.Item(132, 13).Value = -2874.2498
$ws.Cells.Item(132, 14).Value = -33212

$ws.Cells.Item(137, 8).Value = 2185.3333
$ws.Cells.Item(137, 9).Value = 2028
$ws.Cells.Item(137, 11).Value = 6084
$ws.Cells.Item(137, 13).Value = -3534
$ws.Cells.Item(137, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1677
$ws.Cells.Item(2, 9).Value = 142.8
$ws.Cells.Item(2, 11).Value = 142.8
$ws.Cells.Item(2, 13).Value = -29.80000000000001
$ws.Cells.Item(2, 14).ClearContents()

$ws.Cells.Item(5, 8).Value = 408.93332
$ws.Cells.Item(5, 9).Value = 399.1
$ws.Cells.Item(5, 11).Value = 399.1
$ws.Cells.Item(5, 13).Value = -287.1
$ws.Cells.Item(5, 14).ClearContents()

$ws.Cells.Item(30, 8).Value = 2426
$ws.Cells.Item(30, 9).Value = 1504.5
$ws.Cells.Item(30, 11).Value = 1504.5
$ws.Cells.Item(30, 13).Value = -1354.5
$ws.Cells.Item(30, 14).ClearContents()

$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(35, 13).ClearContents()

$ws.Cells.Item(36, 8).Value = 15998.333
$ws.Cells.Item(36, 9).Value = 15998.333
$ws.Cells.Item(36, 11).Value = 15998.333
$ws.Cells.Item(36, 13).Value = -15652.333

$ws.Cells.Item(37, 8).Value = 23750
$ws.Cells.Item(37, 9).Value = 23750
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 11).Value = 23750
$ws.Cells.Item(37, 12).Value = 0
$ws.Cells.Item(37, 13).Value = -23477
$ws.Cells.Item(37, 14).ClearContents()

$ws.Cells.Item(45, 8).Value = 4925.8
$ws.Cells.Item(45, 9).Value = 1814.5
$ws.Cells.Item(45, 11).Value = 1814.5
$ws.Cells.Item(45, 13).Value = -1437.5
$ws.Cells.Item(45, 14).ClearContents()

$ws.Cells.Item(102, 8).Value = 2874.25
$ws.Cells.Item(102, 9).Value = 1499.1666
$ws.Cells.Item(102, 11).Value = 1499.1666
$ws.Cells.Item(102, 13).Value = 122.8334
$ws.Cells.Item(102, 14).ClearContents()

$ws.Cells.Item(116, 8).Value = 1677
$ws.Cells.Item(116, 9).Value = 142.8
$ws.Cells.Item(116, 11).Value = 142.8
$ws.Cells.Item(116, 13).Value = 2151.2
$ws.Cells.Item(116, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1677
$ws.Cells.Item(3, 9).Value = 142.8
$ws.Cells.Item(3, 11).Value = 142.8
$ws.Cells.Item(3, 13).Value = -28.80000000000001
$ws.Cells.Item(3, 14).ClearContents()

$ws.Cells.Item(4, 8).Value = 408.93332
$ws.Cells.Item(4, 9).Value = 399.1
$ws.Cells.Item(4, 11).Value = 399.1
$ws.Cells.Item(4, 13).Value = -284.1
$ws.Cells.Item(4, 14).ClearContents()

$ws.Cells.Item(11, 8).Value = 1049.75
$ws.Cells.Item(11, 9).Value = 100
$ws.Cells.Item(11, 10).Value = 1999.5
$ws.Cells.Item(11, 11).Value = 100
$ws.Cells.Item(11, 12).Value = 1999.5
$ws.Cells.Item(11, 13).Value = 40
$ws.Cells.Item(11, 14).Value = -2279.5

$ws.Cells.Item(14, 8).Value = 10000
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 10).Value = 10000
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 12).Value = 10000
$ws.Cells.Item(14, 13).ClearContents()
$ws.Cells.Item(14, 14).Value = -10344

$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 13).ClearContents()

$ws.Cells.Item(20, 8).Value = 2082.8572
$ws.Cells.Item(20, 9).Value = 2082.8572
$ws.Cells.Item(20, 11).Value = 2082.8572
$ws.Cells.Item(20, 13).Value = -1835.8572

$ws.Cells.Item(22, 8).Value = 572.1
$ws.Cells.Item(22, 9).Value = 403.5
$ws.Cells.Item(22, 10).Value = 825
$ws.Cells.Item(22, 11).Value = 403.5
$ws.Cells.Item(22, 12).Value = 825
$ws.Cells.Item(22, 13).Value = -230.5
$ws.Cells.Item(22, 14).Value = -1171

$ws.Cells.Item(29, 8).Value = 402.5
$ws.Cells.Item(29, 9).Value = 402.5
$ws.Cells.Item(29, 11).Value = 402.5
$ws.Cells.Item(29, 13).Value = -113.5

$ws.Cells.Item(36, 8).Value = 2600
$ws.Cells.Item(36, 9).Value = 2600
$ws.Cells.Item(36, 11).Value = 2600
$ws.Cells.Item(36, 13).Value = -2066

$ws.Cells.Item(37, 8).Value = 701.5
$ws.Cells.Item(37, 9).Value = 701.5
$ws.Cells.Item(37, 11).Value = 701.5
$ws.Cells.Item(37, 13).Value = -564.5

$ws.Cells.Item(80, 8).Value = 512.0833
$ws.Cells.Item(80, 9).Value = 600.5
$ws.Cells.Item(80, 10).Value = 423.66666
$ws.Cells.Item(80, 11).Value = 600.5
$ws.Cells.Item(80, 12).Value = 423.66666
$ws.Cells.Item(80, 13).Value = 397.5
$ws.Cells.Item(80, 14).Value = -2419.66666

$ws.Cells.Item(83, 8).Value = 512.0833
$ws.Cells.Item(83, 9).Value = 600.5
$ws.Cells.Item(83, 10).Value = 423.66666
$ws.Cells.Item(83, 11).Value = 3002.5
$ws.Cells.Item(83, 12).Value = 2118.3333
$ws.Cells.Item(83, 13).Value = 1989.5
$ws.Cells.Item(83, 14).Value = -12102.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 2327.4443
$ws.Cells.Item(7, 10).Value = 2142.6
$ws.Cells.Item(7, 12).Value = 2142.6
$ws.Cells.Item(7, 14).Value = -2368.6

$ws.Cells.Item(22, 8).Value = 1254.2307
$ws.Cells.Item(22, 9).Value = 203.25
$ws.Cells.Item(22, 11).Value = 203.25
$ws.Cells.Item(22, 13).Value = 146.75
$ws.Cells.Item(22, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).ClearContents()
$ws.Cells.Item(16, 14).ClearContents()

$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 14).ClearContents()

$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 14).ClearContents()

$ws.Cells.Item(87, 8).Value = 1000
$ws.Cells.Item(87, 9).Value = 1000
$ws.Cells.Item(87, 11).Value = 3000
$ws.Cells.Item(87, 13).Value = -1752

$ws.Cells.Item(90, 8).Value = 1000
$ws.Cells.Item(90, 9).Value = 1000
$ws.Cells.Item(90, 11).Value = 9000
$ws.Cells.Item(90, 13).Value = -2760

$ws.Cells.Item(97, 8).Value = 981.6667
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 14).ClearContents()

$ws.Cells.Item(112, 8).Value = 3833
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 14).ClearContents()

$ws.Cells.Item(138, 8).Value = 1763.3334
$ws.Cells.Item(138, 9).Value = 1644.5
$ws.Cells.Item(138, 11).Value = 4933.5
$ws.Cells.Item(138, 13).Value = 206.5
$ws.Cells.Item(138, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(13, 8).Value = 226.66667
$ws.Cells.Item(13, 10).Value = 333
$ws.Cells.Item(13, 12).Value = 333
$ws.Cells.Item(13, 14).Value = -611

$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 0
$ws.Cells.Item(17, 13).ClearContents()
$ws.Cells.Item(17, 14).ClearContents()

$ws.Cells.Item(19, 8).Value = 5002500
$ws.Cells.Item(19, 9).Value = 5002500
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 5002500
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = -5002212
$ws.Cells.Item(19, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4996.3335
$ws.Cells.Item(7, 10).Value = 4992
$ws.Cells.Item(7, 12).Value = 4992
$ws.Cells.Item(7, 14).Value = -5216

$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 14).ClearContents()

$ws.Cells.Item(55, 8).Value = 873
$ws.Cells.Item(55, 10).Value = 873
$ws.Cells.Item(55, 12).Value = 873
$ws.Cells.Item(55, 14).Value = -1219

$ws.Cells.Item(126, 8).Value = 4996.3335
$ws.Cells.Item(126, 10).Value = 4992
$ws.Cells.Item(126, 12).Value = 14976
$ws.Cells.Item(126, 14).Value = -19916

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(47, 8).Value = 5000
$ws.Cells.Item(47, 10).Value = 5000
$ws.Cells.Item(47, 12).Value = 5000
$ws.Cells.Item(47, 14).Value = -6144

$ws.Cells.Item(81, 8).Value = 950
$ws.Cells.Item(81, 9).Value = 950
$ws.Cells.Item(81, 11).Value = 1900
$ws.Cells.Item(81, 13).Value = -839

$ws.Cells.Item(84, 8).Value = 950
$ws.Cells.Item(84, 9).Value = 950
$ws.Cells.Item(84, 11).Value = 9500
$ws.Cells.Item(84, 13).Value = -4196

$ws.Cells.Item(113, 8).Value = 675.1667
$ws.Cells.Item(113, 9).Value = 554.7273
$ws.Cells.Item(113, 11).Value = 1664.1819
$ws.Cells.Item(113, 13).Value = 505.8181
$ws.Cells.Item(113, 14).ClearContents()

Write-Host "Applied market data updates."
